# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 3551
    4  = 150
    5  = 7041
    6  = 3699
    7  = 64
    12 = 52
    13 = 38
    16 = 58
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
